$d = $word.ActiveDocument

# Locate the run containing "En esta sección " (trailing space included)
$rng = $d.Content
$rng.Find.Execute("En esta sección ")
$start = $rng.Start
$end = $rng.End

# Insert ", " right after "En esta sección " (i.e. before "se procederá...")
$rng.Collapse(0)
$rng.InsertBefore(", ")

# Remove the original trailing space (the one that was right after "sección")
$spaceRng = $d.Range($end - 1, $end)
$spaceRng.Text = ""

# At this point the paragraph text reads "En esta sección, se procederá..."
# all living inside a single run (text edits make this engine coalesce
# runs that share identical formatting). Force a split into the desired
# three runs -- "En esta sección" | "," | " " -- by toggling Bold on/off
# on the precise sub-ranges; a formatting-only operation splits runs
# without altering the visible formatting, since it's reset right away.
$commaRng = $d.Range($start + 15, $start + 16)
$commaRng.Bold = 1
$commaRng.Bold = 0

$spaceRng2 = $d.Range($start + 16, $start + 17)
$spaceRng2.Bold = 1
$spaceRng2.Bold = 0

# The text edit above also coalesced the two runs that originally
# followed our target run ("se procederá a continuar" and " con el
# desarrollo de la aplicación de tipo contador: ") into a single run.
# Restore that original split point the same way, toggling formatting
# on the whole remainder so the break falls cleanly between the two
# pieces instead of carving out a stray one-character run.
$secondRunLength = ("se procederá a continuar").Length
$boundary = $start + 17 + $secondRunLength
$thirdRunLength = (" con el desarrollo de la aplicación de tipo contador: ").Length
$afterRng = $d.Range($boundary, $boundary + $thirdRunLength)
$afterRng.Bold = 1
$afterRng.Bold = 0
